# Scheduled market-data refresh: update Leve profit figures across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 351211.72
$ws.Range("I15").Value = 351211.72
$ws.Range("K15").Value = 1053635.16
$ws.Range("M15").Value = -1053466.16

$ws.Range("H87").Value = 57354
$ws.Range("J87").Value = 57354
$ws.Range("L87").Value = 57354
$ws.Range("N87").Value = -59850

$ws.Range("H90").Value = 57354
$ws.Range("J90").Value = 57354
$ws.Range("L90").Value = 172062
$ws.Range("N90").Value = -184542

$ws.Range("H112").Value = 1314
$ws.Range("J112").Value = 1357.4468
$ws.Range("L112").Value = 4072.3404
$ws.Range("N112").Value = -6288.3404

$ws.Range("H125").Value = 1560.8889
$ws.Range("I125").Value = 2878
$ws.Range("J125").Value = 507.2
$ws.Range("K125").Value = 25902
$ws.Range("L125").Value = 4564.8
$ws.Range("M125").Value = -23442
$ws.Range("N125").Value = -9484.799999999999

$ws.Range("H134").Value = 58057.332
$ws.Range("J134").Value = 58057.332
$ws.Range("L134").Value = 58057.332
$ws.Range("N134").Value = -68197.33199999999

$ws.Range("H138").Value = 9261167
$ws.Range("I138").Value = 12822040
$ws.Range("J138").Value = 2895
$ws.Range("K138").Value = 38466120
$ws.Range("L138").Value = 8685
$ws.Range("M138").Value = -38460980
$ws.Range("N138").Value = -18965

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9376.65
$ws.Range("I32").Value = 9742.5
$ws.Range("K32").Value = 9742.5
$ws.Range("M32").Value = -9455.5

$ws.Range("H45").Value = 1796.6666
$ws.Range("I45").Value = 1821.4546
$ws.Range("J45").Value = 1728.5
$ws.Range("K45").Value = 1821.4546
$ws.Range("L45").Value = 1728.5
$ws.Range("M45").Value = -1444.4546
$ws.Range("N45").Value = -2482.5

$ws.Range("H61").Value = 8198020.5
$ws.Range("I61").Value = 9435265
$ws.Range("J61").Value = 1276.875
$ws.Range("K61").Value = 9435265
$ws.Range("L61").Value = 1276.875
$ws.Range("M61").Value = -9435053
$ws.Range("N61").Value = -1700.875

$ws.Range("H74").Value = 8335310.5
$ws.Range("I74").Value = 10639717
$ws.Range("J74").Value = 3994.4614
$ws.Range("K74").Value = 10639717
$ws.Range("L74").Value = 3994.4614
$ws.Range("M74").Value = -10638843
$ws.Range("N74").Value = -5742.4614

$ws.Range("H77").Value = 8335310.5
$ws.Range("I77").Value = 10639717
$ws.Range("J77").Value = 3994.4614
$ws.Range("K77").Value = 53198585
$ws.Range("L77").Value = 19972.307
$ws.Range("M77").Value = -53194217
$ws.Range("N77").Value = -28708.307

$ws.Range("H122").Value = 11438.917
$ws.Range("I122").Value = 17187.428
$ws.Range("K122").Value = 51562.284
$ws.Range("M122").Value = -49112.284

$ws.Range("H136").Value = 8198020.5
$ws.Range("I136").Value = 9435265
$ws.Range("J136").Value = 1276.875
$ws.Range("K136").Value = 28305795
$ws.Range("L136").Value = 3830.625
$ws.Range("M136").Value = -28303245
$ws.Range("N136").Value = -8930.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 39300
$ws.Range("J35").Value = 39300
$ws.Range("L35").Value = 39300
$ws.Range("N35").Value = -39920

$ws.Range("H57").Value = 86323.336
$ws.Range("J57").Value = 86323.336
$ws.Range("L57").Value = 86323.336
$ws.Range("N57").Value = -87763.336

$ws.Range("H132").Value = 56929.918
$ws.Range("J132").Value = 56929.918
$ws.Range("L132").Value = 56929.918
$ws.Range("N132").Value = -67049.91800000001

$ws.Range("H134").Value = 2180.9194
$ws.Range("I134").Value = 1160.0444
$ws.Range("J134").Value = 4883.2354
$ws.Range("K134").Value = 3480.1332
$ws.Range("L134").Value = 14649.7062
$ws.Range("M134").Value = -945.1332000000002
$ws.Range("N134").Value = -19719.7062

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H136").Value = 86323.336
$ws.Range("J136").Value = 86323.336
$ws.Range("L136").Value = 86323.336
$ws.Range("N136").Value = -96523.336

$ws.Range("H137").Value = 100000
$ws.Range("J137").Value = 100000
$ws.Range("L137").Value = 100000
$ws.Range("N137").Value = -110200

$ws.Range("H138").Value = 17575.666
$ws.Range("I138").Value = 17575.666
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 17575.666
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -12435.666
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8338742
$ws.Range("I31").Value = 6810.75
$ws.Range("J31").Value = 20836638
$ws.Range("K31").Value = 6810.75
$ws.Range("L31").Value = 20836638
$ws.Range("M31").Value = -6515.75
$ws.Range("N31").Value = -20837228

$ws.Range("H34").Value = 8338742
$ws.Range("I34").Value = 6810.75
$ws.Range("J34").Value = 20836638
$ws.Range("K34").Value = 6810.75
$ws.Range("L34").Value = 20836638
$ws.Range("M34").Value = -6608.75
$ws.Range("N34").Value = -20837042

$ws.Range("H58").Value = 1563.138
$ws.Range("I58").Value = 714.125
$ws.Range("J58").Value = 2608.077
$ws.Range("K58").Value = 714.125
$ws.Range("L58").Value = 2608.077
$ws.Range("M58").Value = -511.125
$ws.Range("N58").Value = -3014.077

$ws.Range("H99").Value = 1679.6
$ws.Range("I99").Value = 1679.6
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1679.6
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -181.5999999999999
$ws.Range("N99").ClearContents()

$ws.Range("H122").Value = 2302.6
$ws.Range("I122").Value = 2302.6
$ws.Range("K122").Value = 6907.799999999999
$ws.Range("M122").Value = -4457.799999999999

$ws.Range("H126").Value = 1679.6
$ws.Range("I126").Value = 1679.6
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5038.799999999999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2568.799999999999
$ws.Range("N126").ClearContents()

$ws.Range("H136").Value = 1563.138
$ws.Range("I136").Value = 714.125
$ws.Range("J136").Value = 2608.077
$ws.Range("K136").Value = 2142.375
$ws.Range("L136").Value = 7824.231000000001
$ws.Range("M136").Value = 407.625
$ws.Range("N136").Value = -12924.231

$ws.Range("H138").Value = 85000
$ws.Range("J138").Value = 85000
$ws.Range("L138").Value = 85000
$ws.Range("N138").Value = -95280

$ws.Range("H140").Value = 39525
$ws.Range("J140").Value = 39525
$ws.Range("L140").Value = 39525
$ws.Range("N140").Value = -49885

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 43.94737
$ws.Range("J12").Value = 46.916668
$ws.Range("L12").Value = 140.750004
$ws.Range("N12").Value = -486.750004

$ws.Range("H23").Value = 92.80768999999999
$ws.Range("I23").Value = 28.1
$ws.Range("J23").Value = 133.25
$ws.Range("K23").Value = 84.30000000000001
$ws.Range("L23").Value = 399.75
$ws.Range("M23").Value = 150.7
$ws.Range("N23").Value = -869.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4137.826
$ws.Range("I102").Value = 4384.3
$ws.Range("J102").Value = 2494.6667
$ws.Range("K102").Value = 4384.3
$ws.Range("L102").Value = 2494.6667
$ws.Range("M102").Value = -2762.3
$ws.Range("N102").Value = -5738.6667

$ws.Range("H133").Value = 70000
$ws.Range("J133").Value = 70000
$ws.Range("L133").Value = 70000
$ws.Range("N133").Value = -80120

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H137").Value = 68853.336
$ws.Range("J137").Value = 68853.336
$ws.Range("L137").Value = 68853.336
$ws.Range("N137").Value = -79053.336

$ws.Range("H138").Value = 59650
$ws.Range("J138").Value = 59650
$ws.Range("L138").Value = 59650
$ws.Range("N138").Value = -69930

$ws.Range("H140").Value = 69996
$ws.Range("J140").Value = 69996
$ws.Range("L140").Value = 69996
$ws.Range("N140").Value = -80356

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 48571.875
$ws.Range("J127").Value = 48571.875
$ws.Range("L127").Value = 48571.875
$ws.Range("N127").Value = -58491.875

$ws.Range("H134").Value = 65429
$ws.Range("J134").Value = 65429
$ws.Range("L134").Value = 65429
$ws.Range("N134").Value = -75569

$ws.Range("H135").Value = 99964.5
$ws.Range("J135").Value = 99964.5
$ws.Range("L135").Value = 99964.5
$ws.Range("N135").Value = -110104.5

$ws.Range("H136").Value = 29419734
$ws.Range("I136").Value = 35716532
$ws.Range("J136").Value = 34668.332
$ws.Range("K136").Value = 107149596
$ws.Range("L136").Value = 104004.996
$ws.Range("M136").Value = -107147046
$ws.Range("N136").Value = -109104.996

$ws.Range("H137").Value = 60412
$ws.Range("J137").Value = 60412
$ws.Range("L137").Value = 60412
$ws.Range("N137").Value = -70612

$ws.Range("H139").Value = 57960
$ws.Range("J139").Value = 57960
$ws.Range("L139").Value = 57960
$ws.Range("N139").Value = -68240

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2493.1667
$ws.Range("I126").Value = 1569.05
$ws.Range("J126").Value = 7113.75
$ws.Range("K126").Value = 4707.15
$ws.Range("L126").Value = 21341.25
$ws.Range("M126").Value = -2237.15
$ws.Range("N126").Value = -26281.25

$ws.Range("H136").Value = 988.02325
$ws.Range("I136").Value = 837.35297
$ws.Range("J136").Value = 1557.2222
$ws.Range("K136").Value = 2512.05891
$ws.Range("L136").Value = 4671.6666
$ws.Range("M136").Value = 37.9410899999998
$ws.Range("N136").Value = -9771.6666
